# Update "Log Loss" values (column B) on the active worksheet to reflect
# the first DoE ML commit results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.11073
$ws.Range("B3").Value = 0.30114
$ws.Range("B4").Value = 0.3012
$ws.Range("B5").Value = 0.3013
$ws.Range("B6").Value = 0.3013
$ws.Range("B7").Value = 0.3013
$ws.Range("B8").Value = 0.30136
$ws.Range("B9").Value = 0.30136
$ws.Range("B10").Value = 0.30141
$ws.Range("B11").Value = 0.30141
$ws.Range("B12").Value = 0.30144
$ws.Range("B13").Value = 0.30152
$ws.Range("B14").Value = 0.30166
$ws.Range("B15").Value = 1.17183
$ws.Range("B16").Value = 1.17212
$ws.Range("B17").Value = 1.17319
$ws.Range("B18").Value = 1.34651
$ws.Range("B19").Value = 2.09828
$ws.Range("B20").Value = 2.25518
$ws.Range("B21").Value = 2.25555
$ws.Range("B22").Value = 2.25555
$ws.Range("B23").Value = 2.25592
$ws.Range("B24").Value = 2.25602
$ws.Range("B25").Value = 2.25639
$ws.Range("B26").Value = 2.25639
$ws.Range("B27").Value = 2.47615
$ws.Range("B28").Value = 4.80288
$ws.Range("B29").Value = 4.80309
$ws.Range("B30").Value = 4.80424
$ws.Range("B31").Value = 4.80486
